# ---------------------------------------------------------------------------
# Applies the "new sim results and new calculation" commit:
#   1. Inserts a brand-new worksheet "sharpe_period" right before "VaR"
#      (so the tab order becomes annualised_return, mean_period_return,
#      sharpe_annualized, sharpe_period, VaR) and fills it with a 10x10
#      correlation-style matrix (same layout/labels as the other sheets).
#   2. Overwrites a handful of recalculated cells on the four pre-existing
#      sheets (annualised_return, mean_period_return, sharpe_annualized,
#      VaR) with their new simulation values.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$entityNames = @(
    "minvar_ports_equalw",
    "minvar_w_cryptos_ahc_equalw",
    "minvar_w_cryptos_kmeans_equalw",
    "minvar_w_cryptos_kshape_equalw",
    "minvar_w_cryptos_random_equalw",
    "rand_ports_equalw",
    "random_w_cryptos_ahc_equalw",
    "random_w_cryptos_kmeans_equalw",
    "random_w_cryptos_kshape_equalw",
    "random_w_cryptos_random_equalw"
)

# ---------------------------------------------------------------------------
# 1. Insert the new "sharpe_period" worksheet before "VaR"
# ---------------------------------------------------------------------------
$varSheet = $wb.Worksheets.Item("VaR")
$newSheet = $wb.Worksheets.Add($varSheet)
$newSheet.Name = "sharpe_period"

# Header row (B1:K1) — bold, centered, boxed, same style as every other sheet
$hdr = $newSheet.Range("B1:K1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
for ($i = 0; $i -lt 10; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $entityNames[$i]
}

# Column A (A2:A11) — same styling, row labels
$lbl = $newSheet.Range("A2:A11")
$lbl.Font.Bold = $true
$lbl.HorizontalAlignment = -4108
$lbl.VerticalAlignment = -4160
$lbl.Borders.LineStyle = 1
for ($i = 0; $i -lt 10; $i++) {
    $newSheet.Cells.Item($i + 2, 1).Value = $entityNames[$i]
}

# 10x10 matrix body (B2:K11)
$matrix = @(
    @(1,      0,      0,      0,      0,      0,      0,      0,      0,      0),
    @(0,      1,      0,      0,      1,      0,      0,      0,      0,      0),
    @(0,      0,      1,      1,      0,      0,      0,      0,      0,      0),
    @(0,      0,      1,      1,      0,      0,      0,      0,      0,      0),
    @(0,      1,      0,      0,      1,      0,      0,      0,      0,      0),
    @(0,      0,      0,      0,      0,      1,      0,      0,      0,      0),
    @(0,      0,      0,      0,      0,      0,      1,      0.0178, 0.1645, 0),
    @(0,      0,      0,      0,      0,      0,      0.0178, 1,      1,      0),
    @(0,      0,      0,      0,      0,      0,      0.1645, 1,      1,      0),
    @(0,      0,      0,      0,      0,      0,      0,      0,      0,      1)
)
for ($r = 0; $r -lt 10; $r++) {
    for ($c = 0; $c -lt 10; $c++) {
        $newSheet.Cells.Item($r + 2, $c + 2).Value = $matrix[$r][$c]
    }
}

# ---------------------------------------------------------------------------
# 2. Update recalculated cells on the existing sheets
# ---------------------------------------------------------------------------

# annualised_return & mean_period_return share the exact same updates
$sameEdits = @{
    "F2"  = 1;
    "F3"  = 0.0325;
    "B6"  = 1;
    "C6"  = 0.0325;
    "H7"  = 0.0183;
    "K7"  = 0.0023;
    "G8"  = 0.0183;
    "J9"  = 0.0131;
    "I10" = 0.0131;
    "G11" = 0.0023
}
foreach ($sheetName in @("annualised_return", "mean_period_return")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $sameEdits.Keys) {
        $ws.Range($addr).Value = $sameEdits[$addr]
    }
}

# sharpe_annualized
$ws3 = $wb.Worksheets.Item("sharpe_annualized")
$edits3 = @{
    "G2"  = 1;
    "E3"  = 0.0001;
    "F3"  = 0.3847;
    "F4"  = 1;
    "K4"  = 0;
    "C5"  = 0.0001;
    "F5"  = 1;
    "K5"  = 0;
    "C6"  = 0.3847;
    "D6"  = 1;
    "E6"  = 1;
    "B7"  = 1;
    "H7"  = 0;
    "I7"  = 0;
    "G8"  = 0;
    "J8"  = 0.0018;
    "K8"  = 0.0004;
    "G9"  = 0;
    "J9"  = 0.0007;
    "K9"  = 0.0011;
    "H10" = 0.0018;
    "I10" = 0.0007;
    "D11" = 0;
    "E11" = 0;
    "H11" = 0.0004;
    "I11" = 0.0011
}
foreach ($addr in $edits3.Keys) {
    $ws3.Range($addr).Value = $edits3[$addr]
}

# VaR
$ws4 = $wb.Worksheets.Item("VaR")
$edits4 = @{
    "G2"  = 1;
    "D3"  = 0.0002;
    "E3"  = 0.0003;
    "F3"  = 0.4958;
    "C4"  = 0.0002;
    "F4"  = 1;
    "G4"  = 0;
    "C5"  = 0.0003;
    "F5"  = 1;
    "G5"  = 0;
    "C6"  = 0.4958;
    "D6"  = 1;
    "E6"  = 1;
    "B7"  = 1;
    "D7"  = 0;
    "E7"  = 0;
    "J8"  = 0.0035;
    "K8"  = 0.0011;
    "J9"  = 0.0008;
    "K9"  = 0.0046;
    "H10" = 0.0035;
    "I10" = 0.0008;
    "K10" = 0;
    "H11" = 0.0011;
    "I11" = 0.0046;
    "J11" = 0
}
foreach ($addr in $edits4.Keys) {
    $ws4.Range($addr).Value = $edits4[$addr]
}

Write-Host "Edit complete"
